# Apply the cell updates described by the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.141.15'
$ws.Range('E2').Value = '  -0.22%  '

$ws.Range('D3').Value = '2.431.19'
$ws.Range('E3').Value = '  +0.38%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = "'563.29"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.13%  '

$ws.Range('D6').Value = "'144.49"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.43%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').Value = '2.431.07'
$ws.Range('E9').Value = '  +0.52%  '

$ws.Range('D10').Value = "'0.110"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.01%  '

$ws.Range('E11').Value = '  +0.23%  '

$ws.Range('D12').Value = "'5.23"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.60%  '

$ws.Range('E13').Value = '  -1.11%  '

$ws.Range('D14').Value = "'26.54"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.58%  '

$ws.Range('E15').Value = '  -1.50%  '

$ws.Range('D16').Value = '2.861.69'
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('D17').Value = '62.086.19'
$ws.Range('E17').Value = '  +0.10%  '

$ws.Range('D18').Value = '2.428.71'
$ws.Range('E18').Value = '  +0.31%  '

$ws.Range('D19').Value = "'11.27"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.83%  '

$ws.Range('D20').Value = "'323.59"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.27%  '

$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = "'4.16"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.01%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'6.84"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.53%  '

$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('D24').Value = "'67.54"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.72%  '

$ws.Range('D25').Value = "'1.72"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.74%  '

$ws.Range('D26').Value = "'8.65"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.58%  '

$ws.Range('D27').Value = "'554.86"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.43%  '

$ws.Range('D28').Value = '2.550.53'
$ws.Range('E28').Value = '  +0.48%  '

$ws.Range('D29').Value = "'0.997"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.34%  '

$ws.Range('D30').Value = '0.0₃0938'
$ws.Range('E30').Value = '  -0.90%  '

$ws.Range('E31').Value = '  +0.08%  '

$ws.Range('D32').Value = "'1.41"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.67%  '

$ws.Range('E33').Value = '  -1.91%  '

$ws.Range('E34').Value = '  -1.27%  '

$ws.Range('E35').Value = '  -2.01%  '

$ws.Range('E36').Value = '  -0.05%  '

$ws.Range('E37').Value = '  -0.14%  '

$ws.Range('D38').Value = "'0.382"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.65%  '

$ws.Range('D39').Value = "'5.54"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.37%  '

$ws.Range('D40').Value = "'151.86"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.10%  '

$ws.Range('D41').Value = "'18.75"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.50%  '

$ws.Range('E42').Value = '  -1.25%  '

$ws.Range('E43').Value = '  +0.20%  '

$ws.Range('D44').Value = "'2.27"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.51%  '

$ws.Range('D45').Value = "'148.18"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.37%  '

$ws.Range('D46').Value = "'3.67"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.03%  '

$ws.Range('D47').Value = "'0.0532"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.99%  '

$ws.Range('D48').Value = "'20.11"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.75%  '

$ws.Range('D49').Value = "'0.597"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.39%  '

$ws.Range('D50').Value = "'0.0925"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.28%  '

$ws.Range('D51').Value = "'0.0230"
$ws.Range('D51').ClearFormats()

